# Nayem_meal.xlsx edit: tel (oil) 500g -98, Minhaz raat/vat short -> meal -1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Minhaz (row 6) - raat/vat (dinner) was short: V6 meal count 2.5 -> 1.5
$ws.Range("V6").Value = 1.5

# Bazar (shopping) row 42: Mahfuz bought something on the "V" day column
$ws.Range("V42").Value = "Mahfuz"

# Bazar TK (row 43): tel (cooking oil) 500g cost 98 taka, added to that day's bazar cost
$ws.Range("V43").Value = 98

# Reflect the edited cell as the active selection, matching the live edit session
$ws.Range("V6").Select()
